$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (humoment -> shape)
$ws.Range("J1").Value = "shape"
$ws.Range("R1").Value = "texture-shape"
$ws.Range("V1").Value = "color-shape"
$ws.Range("Z1").Value = "texture-color-shape"

# Update numeric values for rows 4-7

# Row 4
$ws.Range("B4").Value = 0.5932486733548871
$ws.Range("C4").Value = 0.622
$ws.Range("D4").Value = 0.580538685613357
$ws.Range("E4").Value = 0.604
$ws.Range("F4").Value = 0.6084982694032386
$ws.Range("G4").Value = 0.6499999999999999
$ws.Range("H4").Value = 0.6151592825421899
$ws.Range("I4").Value = 0.6165
$ws.Range("J4").Value = 0.6698229371633928
$ws.Range("K4").Value = 0.9799999999999999
$ws.Range("L4").Value = 0.5089641738515409
$ws.Range("M4").Value = 0.517
$ws.Range("N4").Value = 0.6436996983212654
$ws.Range("O4").Value = 0.6599999999999999
$ws.Range("P4").Value = 0.6641647457149111
$ws.Range("Q4").Value = 0.6570000000000001
$ws.Range("R4").Value = 0.5900637650000664
$ws.Range("S4").Value = 0.6199999999999999
$ws.Range("T4").Value = 0.57565129568173
$ws.Range("U4").Value = 0.6004999999999999
$ws.Range("V4").Value = 0.601833407479299
$ws.Range("W4").Value = 0.6449999999999999
$ws.Range("X4").Value = 0.6072180444423341
$ws.Range("Y4").Value = 0.61
$ws.Range("Z4").Value = 0.6420919044424304
$ws.Range("AA4").Value = 0.659
$ws.Range("AB4").Value = 0.6623320251403092
$ws.Range("AC4").Value = 0.655

# Row 5
$ws.Range("B5").Value = 0.5963164773123009
$ws.Range("C5").Value = 0.624
$ws.Range("D5").Value = 0.5825920123720698
$ws.Range("E5").Value = 0.6055
$ws.Range("F5").Value = 0.7035724931051273
$ws.Range("G5").Value = 0.765
$ws.Range("H5").Value = 0.6674766636113154
$ws.Range("I5").Value = 0.675
$ws.Range("J5").Value = 0.6525963824887014
$ws.Range("K5").Value = 0.857
$ws.Range("L5").Value = 0.5320614704031437
$ws.Range("M5").Value = 0.5485
$ws.Range("N5").Value = 0.6389404811190171
$ws.Range("O5").Value = 0.652
$ws.Range("P5").Value = 0.6481075370780268
$ws.Range("Q5").Value = 0.6519999999999999
$ws.Range("R5").Value = 0.5964912338594319
$ws.Range("S5").Value = 0.6239999999999999
$ws.Range("T5").Value = 0.5829737464980567
$ws.Range("U5").Value = 0.6054999999999999
$ws.Range("V5").Value = 0.6832539814881928
$ws.Range("W5").Value = 0.7470000000000001
$ws.Range("X5").Value = 0.6493748164999953
$ws.Range("Y5").Value = 0.6525000000000001
$ws.Range("Z5").Value = 0.6391019338593521
$ws.Range("AA5").Value = 0.6540000000000001
$ws.Range("AB5").Value = 0.6471785098518685
$ws.Range("AC5").Value = 0.652

# Row 6
$ws.Range("B6").Value = 0.6036519292314397
$ws.Range("C6").Value = 0.635
$ws.Range("D6").Value = 0.5912403920193882
$ws.Range("E6").Value = 0.6134999999999999
$ws.Range("F6").Value = 0.7065737934675792
$ws.Range("G6").Value = 0.74
$ws.Range("H6").Value = 0.6930043549100147
$ws.Range("I6").Value = 0.6944999999999999
$ws.Range("J6").Value = 0.6580562848705325
$ws.Range("K6").Value = 0.885
$ws.Range("L6").Value = 0.5272875985295125
$ws.Range("M6").Value = 0.543
$ws.Range("N6").Value = 0.6648943608648168
$ws.Range("O6").Value = 0.6699999999999999
$ws.Range("P6").Value = 0.6811693287624798
$ws.Range("Q6").Value = 0.68
$ws.Range("R6").Value = 0.6056245702751673
$ws.Range("S6").Value = 0.6380000000000001
$ws.Range("T6").Value = 0.5927528505075299
$ws.Range("U6").Value = 0.6165
$ws.Range("V6").Value = 0.6635843691012439
$ws.Range("W6").Value = 0.706
$ws.Range("X6").Value = 0.6646969646848756
$ws.Range("Y6").Value = 0.653
$ws.Range("Z6").Value = 0.655087619539876
$ws.Range("AA6").Value = 0.657
$ws.Range("AB6").Value = 0.6761167703753674
$ws.Range("AC6").Value = 0.673

# Row 7
$ws.Range("B7").Value = 0.4765152655222479
$ws.Range("C7").Value = 0.488
$ws.Range("D7").Value = 0.473380765413255
$ws.Range("E7").Value = 0.4784999999999999
$ws.Range("F7").Value = 0.4839927403138217
$ws.Range("G7").Value = 0.499
$ws.Range("H7").Value = 0.4805376887905241
$ws.Range("I7").Value = 0.4834999999999999
$ws.Range("J7").Value = 0.6183677339517268
$ws.Range("K7").Value = 0.8370000000000001
$ws.Range("L7").Value = 0.4996152807756056
$ws.Range("M7").Value = 0.505
$ws.Range("N7").Value = 0.5178702970339001
$ws.Range("O7").Value = 0.541
$ws.Range("P7").Value = 0.511124862071197
$ws.Range("Q7").Value = 0.5170000000000001
$ws.Range("R7").Value = 0.4907741497664189
$ws.Range("S7").Value = 0.506
$ws.Range("T7").Value = 0.484064698609001
$ws.Range("U7").Value = 0.487
$ws.Range("V7").Value = 0.5347117507663307
$ws.Range("W7").Value = 0.5559999999999999
$ws.Range("X7").Value = 0.5410911672502035
$ws.Range("Y7").Value = 0.5475000000000001
$ws.Range("Z7").Value = 0.5260243126713876
$ws.Range("AA7").Value = 0.549
$ws.Range("AB7").Value = 0.5198667432613931
$ws.Range("AC7").Value = 0.526
